$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reformat the "Related Publication" citation in B21: collapse the
# multi-line citation text into a single line (formatting cleanup).
$newCitation = "Aleva, A., van den Berg, T., Laceulle, O.M. et al. A smartphone-based intervention for young people who self-harm (‘PRIMARY’): study protocol for a multicenter randomized controlled trial. BMC Psychiatry 23, 840 (2023). doi: <https://doi.org/10.1186/s12888-023-05301-x>"
$ws.Range("B21").Value = $newCitation

# The row no longer needs the tall explicit height that accommodated the
# old 3-line wrapped text, so auto-fit it back down to the default.
$ws.Rows.Item(21).EntireRow.AutoFit()

# Update view state: scroll position and current selection.
$ws.Range("B25").Select()

Write-Output "done"
